$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 485, shifting rows 485:525 down to 486:526
$ws.Rows.Item(485).Insert()

$ws.Cells.Item(485, 1).Value = 3
$ws.Cells.Item(485, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(485, 3).Value = "Coquimbo"
$ws.Cells.Item(485, 4).Value = 45013
$ws.Cells.Item(485, 5).Value = 5
$ws.Cells.Item(485, 6).Value = 100114013
$ws.Cells.Item(485, 7).Value = "Zanahoria"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 340
$ws.Cells.Item(485, 11).Value = 9000
$ws.Cells.Item(485, 12).Value = 9500
$ws.Cells.Item(485, 13).Value = 9265
$ws.Cells.Item(485, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(485, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(485, 16).Value = 463
$ws.Cells.Item(485, 17).Value = 20
$ws.Cells.Item(485, 18).Value = "Hortaliza"
